$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "Despesa"
$ws.Range("B7").Value = "SERVIÇOS"
$ws.Range("C7").Value = 50
$ws.Range("D7").Value = "26/01/2025"

$ws.Range("A8").Value = "Receita"
$ws.Range("B8").Value = "ALUGUEL"
$ws.Range("C8").Value = 400
$ws.Range("D8").Value = "26/01/2025"
